$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact literal text value, bypassing Excel's
# automatic number/date inference (values like "1.00" or "67.590.66" must
# stay as text, matching the source data which stores these as strings).
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '67.590.66'
$ws.Range("E2").Value = '  -3.67%  '
Set-TextValue $ws.Range("D3") '3.283.63'
$ws.Range("E3").Value = '  -5.76%  '
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  +0.17%  '
Set-TextValue $ws.Range("D5") '595.38'
$ws.Range("E5").Value = '  -3.27%  '
Set-TextValue $ws.Range("D6") '152.53'
$ws.Range("E6").Value = '  -9.68%  '
$ws.Range("E7").Value = '  -0.15%  '
Set-TextValue $ws.Range("D8") '3.276.89'
$ws.Range("E8").Value = '  -5.85%  '
Set-TextValue $ws.Range("D9") '0.547'
$ws.Range("E9").Value = '  -8.78%  '
Set-TextValue $ws.Range("D10") '0.172'
$ws.Range("E10").Value = '  -11.77%  '
$ws.Range("E11").Value = '  -3.71%  '
Set-TextValue $ws.Range("D12") '0.511'
$ws.Range("E12").Value = '  -10.43%  '
Set-TextValue $ws.Range("D13") '38.88'
$ws.Range("E13").Value = '  -13.48%  '
$ws.Range("E14").Value = '  -8.85%  '
Set-TextValue $ws.Range("D15") '3.804.99'
$ws.Range("E15").Value = '  -5.86%  '
Set-TextValue $ws.Range("D16") '67.610.15'
$ws.Range("E16").Value = '  -3.72%  '
Set-TextValue $ws.Range("D17") '3.282.31'
$ws.Range("E17").Value = '  -5.84%  '
Set-TextValue $ws.Range("D18") '538.49'
$ws.Range("E18").Value = '  -9.23%  '
Set-TextValue $ws.Range("D19") '0.115'
$ws.Range("E19").Value = '  -5.38%  '
$ws.Range("E20").Value = '  -12.62%  '
Set-TextValue $ws.Range("D21") '15.19'
$ws.Range("E21").Value = '  -12.26%  '
Set-TextValue $ws.Range("D22") '0.767'
$ws.Range("E22").Value = '  -11.38%  '
Set-TextValue $ws.Range("D23") '7.92'
$ws.Range("E23").Value = '  -10.52%  '
Set-TextValue $ws.Range("D24") '86.16'
$ws.Range("E24").Value = '  -10.50%  '
Set-TextValue $ws.Range("D25") '13.69'
$ws.Range("E25").Value = '  -10.31%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D26") '3.28'
$ws.Range("E26").Value = '  -10.01%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D27") '1.00'
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D28") '8.16'
$ws.Range("E28").Value = '  -6.68%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D29") '2.18'
$ws.Range("E29").Value = '  -12.78%  '
Set-TextValue $ws.Range("D30") '29.47'
$ws.Range("E30").Value = '  -11.45%  '
Set-TextValue $ws.Range("D31") '2.72'
$ws.Range("E31").Value = '  -4.49%  '
$ws.Range("E32").Value = '  -7.46%  '
Set-TextValue $ws.Range("D33") '6.67'
$ws.Range("E33").Value = '  -16.15%  '
Set-TextValue $ws.Range("D34") '5.83'
$ws.Range("E34").Value = '  -12.33%  '
Set-TextValue $ws.Range("D35") '532.97'
$ws.Range("E35").Value = '  -6.80%  '
$ws.Range("E36").Value = '  -0.03%  '
Set-TextValue $ws.Range("D37") '0.0456'
$ws.Range("E37").Value = '  -7.01%  '
Set-TextValue $ws.Range("D38") '53.62'
$ws.Range("E38").Value = '  -4.81%  '
Set-TextValue $ws.Range("D39") '0.0864'
$ws.Range("E39").Value = '  -11.00%  '
Set-TextValue $ws.Range("D40") '9.09'
$ws.Range("E40").Value = '  -15.31%  '
$ws.Range("E41").Value = '  -10.10%  '
Set-TextValue $ws.Range("D42") '2.85'
$ws.Range("E42").Value = '  -12.34%  '
Set-TextValue $ws.Range("D43") '2.958.19'
$ws.Range("E43").Value = '  -10.08%  '
Set-TextValue $ws.Range("D44") '0.271'
$ws.Range("E44").Value = '  -10.18%  '
Set-TextValue $ws.Range("D45") '0.0₃0598'
$ws.Range("E45").Value = '  -14.96%  '
$ws.Range("E46").Value = '  -8.65%  '
Set-TextValue $ws.Range("D47") '26.97'
$ws.Range("E47").Value = '  -12.96%  '
$ws.Range("E48").Value = '  -14.82%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("E50").Value = '  -9.99%  '
Set-TextValue $ws.Range("D51") '124.57'
$ws.Range("E51").Value = '  -6.44%  '
